# Updated cryptos list on Thu Mar 21 07:56:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.988.96'
$ws.Range('E2').Value = '  +6.53%  '
$ws.Range('D3').Value = '3.514.27'
$ws.Range('E3').Value = '  +9.48%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '188.65'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +8.91%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '550.20'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.45%  '
$ws.Range('D7').Value = '3.511.30'
$ws.Range('E7').Value = '  +9.43%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.605'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.65%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.631'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.54%  '
$ws.Range('E11').Value = '  +14.56%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '54.53'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000268'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.97%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.34'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.47%  '
$ws.Range('D15').Value = '4.085.70'
$ws.Range('E15').Value = '  +9.75%  '
$ws.Range('D16').Value = '3.521.84'
$ws.Range('E16').Value = '  +9.87%  '
$ws.Range('E17').Value = '  +3.04%  '
$ws.Range('D18').Value = '67.075.89'
$ws.Range('E18').Value = '  +7.00%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '18.13'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.54%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.84'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +6.66%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.992'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.31%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '429.90'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +16.93%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '84.88'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +4.33%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.90'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.23%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.15'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +6.81%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.10'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.35%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.89'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +8.35%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.99'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +5.60%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.95'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +8.86%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '30.16'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +5.83%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '643.06'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.65'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.68'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.87%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.110'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.26%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '59.28'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.78%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '38.42'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +3.46%  '
$ws.Range('D37').Value = '0.0₃0810'
$ws.Range('E37').Value = '  +13.78%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.998'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.389'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.97%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.141'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +13.98%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.35'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +14.44%  '
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('D43').Value = '3.022.47'
$ws.Range('E43').Value = '  +4.72%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.63'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.82%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.90'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +9.17%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.86'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +9.53%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.30'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +12.17%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0416'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +5.75%  '
$ws.Range('E49').Value = '  +4.82%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.67'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +12.32%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '140.95'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +5.23%  '
